$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.460.73'
$ws.Range("E2").Value = '  -0.25%  '

$ws.Range("D3").Value = '1.899.66'
$ws.Range("E3").Value = '  +1.34%  '

$ws.Range("E4").Value = '  +0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '238.60'
$ws.Range("E5").Value = '  +0.78%  '

$ws.Range("E6").Value = '  +0.14%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4901'
$ws.Range("E7").Value = '  +0.60%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2919'
$ws.Range("E8").Value = '  +0.89%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06672'
$ws.Range("E9").Value = '  +0.08%  '

$ws.Range("D10").Value = '1.893.29'
$ws.Range("E10").Value = '  +1.03%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '16.91'
$ws.Range("E11").Value = '  +1.80%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07325'
$ws.Range("E12").Value = '  +1.44%  '

$ws.Range("E13").Value = '  +3.55%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '88.11'
$ws.Range("E14").Value = '  -1.55%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6660'
$ws.Range("E15").Value = '  +1.89%  '

$ws.Range("D16").Value = '30.437.61'
$ws.Range("E16").Value = '  -0.15%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000007863'
$ws.Range("E17").Value = '  +0.39%  '

$ws.Range("E18").Value = '  +3.19%  '

$ws.Range("E19").Value = '  +0.16%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.440'
$ws.Range("E20").Value = '  +15.09%  '

$ws.Range("D21").Value = '2.136.11'
$ws.Range("E21").Value = '  +1.07%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9999'
$ws.Range("E22").Value = '  +0.40%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '195.53'
$ws.Range("E23").Value = '  -8.54%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.129'
$ws.Range("E24").Value = '  +0.01%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.484'
$ws.Range("E25").Value = '  +1.27%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '162.92'
$ws.Range("E26").Value = '  +4.40%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.35'
$ws.Range("E27").Value = '  -3.57%  '

$ws.Range("E28").Value = '  +6.05%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.482'
$ws.Range("E29").Value = '  +5.09%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.324'
$ws.Range("E30").Value = '  +1.47%  '

$ws.Range("E31").Value = '  +1.25%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.112'
$ws.Range("E32").Value = '  +4.77%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05160'

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7386'
$ws.Range("E34").Value = '  +2.10%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.105'
$ws.Range("E35").Value = '  +2.50%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.730'
$ws.Range("E36").Value = '  +1.68%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01841'
$ws.Range("E37").Value = '  +1.62%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.675'
$ws.Range("E38").Value = '  +0.69%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.9236'
$ws.Range("E39").Value = '  +0.48%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.067'
$ws.Range("E40").Value = '  +0.97%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4401'
$ws.Range("E41").Value = '  -0.06%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '106.88'
$ws.Range("E42").Value = '  +2.31%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.895'
$ws.Range("E43").Value = '  +2.80%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9950'
$ws.Range("E44").Value = '  +0.08%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '68.89'
$ws.Range("E45").Value = '  +20.78%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1369'
$ws.Range("E46").Value = '  +3.26%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.564'
$ws.Range("E47").Value = '  +3.05%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.025'
$ws.Range("E48").Value = '  +4.95%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '34.92'
$ws.Range("E49").Value = '  +5.31%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05833'
$ws.Range("E50").Value = '  +0.17%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3923'
$ws.Range("E51").Value = '  -2.30%  '
